$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 328.57144
$ws.Range("I6").Value = 328.57144
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 985.71432
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -873.71432
$ws.Range("N6").Value = ""

$ws.Range("H57").Value = 111001
$ws.Range("J57").Value = 111001
$ws.Range("L57").Value = 333003
$ws.Range("N57").Value = -334001

$ws.Range("H80").Value = 8960.083000000001
$ws.Range("I80").Value = 703.3333
$ws.Range("J80").Value = 17216.834
$ws.Range("K80").Value = 2109.9999
$ws.Range("L80").Value = 51650.50199999999
$ws.Range("M80").Value = -1111.9999
$ws.Range("N80").Value = -53646.50199999999

$ws.Range("H83").Value = 8960.083000000001
$ws.Range("I83").Value = 703.3333
$ws.Range("J83").Value = 17216.834
$ws.Range("K83").Value = 6329.9997
$ws.Range("L83").Value = 154951.506
$ws.Range("M83").Value = -1337.9997
$ws.Range("N83").Value = -164935.506

$ws.Range("H125").Value = 705.1818
$ws.Range("I125").Value = 645.1667
$ws.Range("J125").Value = 777.2
$ws.Range("K125").Value = 5806.5003
$ws.Range("L125").Value = 6994.8
$ws.Range("M125").Value = -3346.5003
$ws.Range("N125").Value = -11914.8

$ws.Range("H132").Value = 4763893.5
$ws.Range("I132").Value = 1518.5079
$ws.Range("K132").Value = 4555.5237
$ws.Range("M132").Value = -2025.5237

$ws.Range("H133").Value = 52097.5
$ws.Range("J133").Value = 52097.5
$ws.Range("L133").Value = 52097.5
$ws.Range("N133").Value = -62217.5

$ws.Range("H136").Value = 56697.777
$ws.Range("J136").Value = 56697.777
$ws.Range("L136").Value = 56697.777
$ws.Range("N136").Value = -66897.777

$ws.Range("H140").Value = 66483.336
$ws.Range("J140").Value = 66483.336
$ws.Range("L140").Value = 66483.336
$ws.Range("N140").Value = -76843.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 101
$ws.Range("I4").Value = 101
$ws.Range("K4").Value = 101
$ws.Range("M4").Value = 15

$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").Value = ""

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = ""

$ws.Range("H32").Value = 18055.346
$ws.Range("I32").Value = 14172.017
$ws.Range("K32").Value = 14172.017
$ws.Range("M32").Value = -13885.017

$ws.Range("H110").Value = 8751
$ws.Range("I110").Value = 9855.053
$ws.Range("J110").Value = 3506.75
$ws.Range("K110").Value = 9855.053
$ws.Range("L110").Value = 3506.75
$ws.Range("M110").Value = -7810.053
$ws.Range("N110").Value = -7596.75

$ws.Range("H132").Value = 3790.7715
$ws.Range("I132").Value = 2840.6667
$ws.Range("J132").Value = 4796.7646
$ws.Range("K132").Value = 8522.000100000001
$ws.Range("L132").Value = 14390.2938
$ws.Range("M132").Value = -5992.000100000001
$ws.Range("N132").Value = -19450.2938

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 57897.75
$ws.Range("J74").Value = 57897.75
$ws.Range("L74").Value = 57897.75
$ws.Range("N74").Value = -59769.75

$ws.Range("H77").Value = 57897.75
$ws.Range("J77").Value = 57897.75
$ws.Range("L77").Value = 173693.25
$ws.Range("N77").Value = -183053.25

$ws.Range("H81").Value = 42539.8
$ws.Range("J81").Value = 42539.8
$ws.Range("L81").Value = 42539.8
$ws.Range("N81").Value = -44661.8

$ws.Range("H84").Value = 42539.8
$ws.Range("J84").Value = 42539.8
$ws.Range("L84").Value = 127619.4
$ws.Range("N84").Value = -138227.4

$ws.Range("H86").Value = 2097.8333
$ws.Range("I86").Value = 1997.4
$ws.Range("J86").Value = 2600
$ws.Range("K86").Value = 1997.4
$ws.Range("L86").Value = 2600
$ws.Range("M86").Value = -874.4000000000001
$ws.Range("N86").Value = -4846

$ws.Range("H89").Value = 2097.8333
$ws.Range("I89").Value = 1997.4
$ws.Range("J89").Value = 2600
$ws.Range("K89").Value = 9987
$ws.Range("L89").Value = 13000
$ws.Range("M89").Value = -4371
$ws.Range("N89").Value = -24232

$ws.Range("H134").Value = 20796.23
$ws.Range("I134").Value = 1508.125
$ws.Range("J134").Value = 252253.5
$ws.Range("K134").Value = 4524.375
$ws.Range("L134").Value = 756760.5
$ws.Range("M134").Value = -1989.375
$ws.Range("N134").Value = -761830.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9899.299999999999
$ws.Range("I3").Value = 18374.285
$ws.Range("J3").Value = 5335.846
$ws.Range("K3").Value = 55122.855
$ws.Range("L3").Value = 16007.538
$ws.Range("M3").Value = -55010.855
$ws.Range("N3").Value = -16231.538

$ws.Range("H7").Value = 1280.75
$ws.Range("J7").Value = 2531
$ws.Range("L7").Value = 7593
$ws.Range("N7").Value = -7817

$ws.Range("H22").Value = 1058.25
$ws.Range("J22").Value = 2733
$ws.Range("L22").Value = 8199
$ws.Range("N22").Value = -8537

$ws.Range("H27").Value = 1058.25
$ws.Range("J27").Value = 2733
$ws.Range("L27").Value = 8199
$ws.Range("N27").Value = -8403

$ws.Range("H92").Value = 803.7143
$ws.Range("J92").Value = 891
$ws.Range("L92").Value = 2673
$ws.Range("N92").Value = -5169

$ws.Range("H122").Value = 5992.864
$ws.Range("I122").Value = 1080.6
$ws.Range("J122").Value = 10086.417
$ws.Range("K122").Value = 9725.4
$ws.Range("L122").Value = 90777.753
$ws.Range("M122").Value = -7275.4
$ws.Range("N122").Value = -95677.753

$ws.Range("H131").Value = 2041860.6
$ws.Range("I131").Value = 7143532
$ws.Range("J131").Value = 1191.9429
$ws.Range("K131").Value = 21430596
$ws.Range("L131").Value = 3575.8287
$ws.Range("M131").Value = -21425556
$ws.Range("N131").Value = -13655.8287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 59341430
$ws.Range("I122").Value = 59157896
$ws.Range("J122").Value = 60002160
$ws.Range("K122").Value = 177473688
$ws.Range("L122").Value = 180006480
$ws.Range("M122").Value = -177471238
$ws.Range("N122").Value = -180011380

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3328.9092
$ws.Range("I7").Value = 2013.5
$ws.Range("J7").Value = 6836.6665
$ws.Range("K7").Value = 2013.5
$ws.Range("L7").Value = 6836.6665
$ws.Range("M7").Value = -1901.5
$ws.Range("N7").Value = -7060.6665

$ws.Range("H40").Value = 3000
$ws.Range("I40").Value = 3000
$ws.Range("K40").Value = 3000
$ws.Range("M40").Value = -2864

$ws.Range("H122").Value = 4289461
$ws.Range("I122").Value = 5498597.5
$ws.Range("K122").Value = 16495792.5
$ws.Range("M122").Value = -16493342.5

$ws.Range("H126").Value = 3328.9092
$ws.Range("I126").Value = 2013.5
$ws.Range("J126").Value = 6836.6665
$ws.Range("K126").Value = 6040.5
$ws.Range("L126").Value = 20509.9995
$ws.Range("M126").Value = -3570.5
$ws.Range("N126").Value = -25449.9995

$ws.Range("H132").Value = 3647.5881
$ws.Range("I132").Value = 3500.9062
$ws.Range("K132").Value = 10502.7186
$ws.Range("M132").Value = -7972.7186

$ws.Range("H139").Value = 67269.164
$ws.Range("J139").Value = 67269.164
$ws.Range("L139").Value = 67269.164
$ws.Range("N139").Value = -77549.164

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1152.0667
$ws.Range("I126").Value = 1134.1818
$ws.Range("J126").Value = 1201.25
$ws.Range("K126").Value = 3402.5454
$ws.Range("L126").Value = 3603.75
$ws.Range("M126").Value = -932.5454
$ws.Range("N126").Value = -8543.75

$ws.Range("H132").Value = 2509.24
$ws.Range("I132").Value = 1004.8
$ws.Range("K132").Value = 3014.4
$ws.Range("M132").Value = -484.3999999999996

$ws.Range("H133").Value = 42888.332
$ws.Range("J133").Value = 42888.332
$ws.Range("L133").Value = 42888.332
$ws.Range("N133").Value = -53008.332
